$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 2

$ws.Range("H14").Select()
